# Business microservices Student and School setup — add Port column and
# two new microservice blocks (Student, School) to the
# SchoolManagementProject sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SchoolManagementProject")

# --- Student microservice block (rows 7-13) ---
$ws.Range("A7").Value = "Student"
$ws.Range("B7").Value = "Postgre sql"
$ws.Range("B8").Value = "Lombok"
$ws.Range("B9").Value = "Spring data jpa"
$ws.Range("B10").Value = "Spring Web"
$ws.Range("B11").Value = "config client"
$ws.Range("B12").Value = "Eureka discovery client"
$ws.Range("B13").Value = "Spring Boot Actuator"
$ws.Range("C7").Value = 8091

# --- Header row: add the new "Port" column header (C1), copying the
# existing header formatting (style s="8") from B1. ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "Port"

# --- Existing Gateway block gains a Port value ---
$ws.Range("C2").Value = 8090

# --- School microservice block (rows 15-21) ---
$ws.Range("A15").Value = "School"
$ws.Range("B15").Value = "Postgre sql"
$ws.Range("B16").Value = "Lombok"
$ws.Range("B17").Value = "Spring data jpa"
$ws.Range("B18").Value = "Spring Web"
$ws.Range("B19").Value = "config client"
$ws.Range("B20").Value = "Eureka discovery client"
$ws.Range("B21").Value = "Spring Boot Actuator"
$ws.Range("C15").Value = 8092

# Leave the active selection on C15, matching the authored workbook state.
$ws.Range("C15").Select() | Out-Null
